$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data per the Jan 4 2024 GitHub Actions refresh

$ws.Range('D2').Value = '44.416.86'
$ws.Range('E2').Value = '  +3.63%  '
$ws.Range('D3').Value = '2.273.93'
$ws.Range('E3').Value = '  +3.14%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '322.04'
$ws.Range('E5').Value = '  +2.07%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '105.21'
$ws.Range('E6').Value = '  +6.58%  '
$ws.Range('E7').Value = '  +0.80%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.572'
$ws.Range('E9').Value = '  +2.61%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.70'
$ws.Range('E10').Value = '  +5.36%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0844'
$ws.Range('E11').Value = '  +2.50%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.91'
$ws.Range('E12').Value = '  +3.36%  '
$ws.Range('E13').Value = '  +0.19%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.884'
$ws.Range('D15').Value = '2.623.13'
$ws.Range('E15').Value = '  +3.19%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.58'
$ws.Range('E16').Value = '  +2.89%  '
$ws.Range('D17').Value = '2.273.47'
$ws.Range('E17').Value = '  +3.22%  '
$ws.Range('D18').Value = '44.356.30'
$ws.Range('E18').Value = '  +3.72%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.90'
$ws.Range('E19').Value = '  -3.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000100'
$ws.Range('E20').Value = '  +4.89%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.54'
$ws.Range('E21').Value = '  +2.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '66.51'
$ws.Range('E22').Value = '  +2.21%  '
$ws.Range('E23').Value = '  +1.94%  '
$ws.Range('E24').Value = '  +2.09%  '
$ws.Range('E25').Value = '  +5.91%  '
$ws.Range('E26').Value = '  +0.37%  '
$ws.Range('E27').Value = '  +3.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '38.43'
$ws.Range('E28').Value = '  +12.24%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.21'
$ws.Range('E29').Value = '  -0.55%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.50'
$ws.Range('E30').Value = '  +3.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.72'
$ws.Range('E31').Value = '  +1.22%  '
$ws.Range('B32').Value = 'Monero'
$ws.Range('C32').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '162.33'
$ws.Range('E32').Value = '  +5.42%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0885'
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('E34').Value = '  -0.52%  '
$ws.Range('E35').Value = '  +9.93%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.01'
$ws.Range('E36').Value = '  +5.95%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.14'
$ws.Range('E37').Value = '  +2.90%  '
$ws.Range('E38').Value = '  +0.63%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.94'
$ws.Range('E39').Value = '  +4.92%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.43'
$ws.Range('E40').Value = '  +0.77%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0330'
$ws.Range('E41').Value = '  +1.90%  '
$ws.Range('B42').Value = 'Celestia'
$ws.Range('C42').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '15.57'
$ws.Range('E42').Value = '  +27.54%  '
$ws.Range('E43').Value = '  +0.15%  '
$ws.Range('D44').Value = '1.788.14'
$ws.Range('E44').Value = '  -3.01%  '
$ws.Range('E45').Value = '  +1.55%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '86.42'
$ws.Range('E46').Value = '  -1.65%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.45'
$ws.Range('E47').Value = '  +1.97%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '60.85'
$ws.Range('E48').Value = '  +0.59%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.74'
$ws.Range('E49').Value = '  +10.28%  '
$ws.Range('B50').Value = 'ordi'
$ws.Range('C50').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '75.81'
$ws.Range('E50').Value = '  +1.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '104.38'
$ws.Range('E51').Value = '  +1.90%  '
